$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet holds a single Excel Table ("Table2") with columns:
# Question | Difficulty | Pattern | Notes | Link
$tbl = $ws.ListObjects.Item(1)

# --- Row 40: 739. Daily Temperatures -----------------------------------
$newRow1 = $tbl.ListRows.Add()
$ws.Range("A38:E38").Copy()
$ws.Range("A40:E40").PasteSpecial(-4122)  # xlPasteFormats (copy "Medium" look)
$ws.Range("A40").Value2 = "739. Daily Temperatures"
$ws.Range("B40").Value2 = "Medium"
$ws.Range("C40").Value2 = "Stack"
$ws.Range("D40").Value2 = "A Monotonic Stack question. Store ""slow"" pointer values in stack and search them when ""fast"" pointer fails to meet a condition."
$ws.Hyperlinks.Add($ws.Range("E40"), "https://leetcode.com/problems/daily-temperatures/solutions/109832/java-easy-ac-solution-with-stack/ ")
$ws.Range("E40").Style = "Hyperlink"

# --- Row 41: 567. Permutation in String ---------------------------------
$newRow2 = $tbl.ListRows.Add()
$ws.Range("A38:E38").Copy()
$ws.Range("A41:E41").PasteSpecial(-4122)  # xlPasteFormats (copy "Medium" look)
$ws.Range("A41").Value2 = "567. Permutation in String"
$ws.Range("B41").Value2 = "Medium"
$ws.Range("C41").Value2 = "Sliding Window"
$ws.Range("D41").Value2 = "Maintain 2 arrays (or hashmaps) for char frequencies, 1 of the target, and 1 of the window. Update the frequencies of the window and check if they match at each iteration."
$ws.Hyperlinks.Add($ws.Range("E41"), "https://leetcode.com/problems/permutation-in-string/solutions/102588/java-solution-sliding-window/ ")
$ws.Range("E41").Style = "Hyperlink"

$ws.Application.CutCopyMode = $false

# --- Selection / scroll position matches the authored edit --------------
$ws.Range("D48").Select()
